$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for the first md file, shared between the
# Overview sheet and the de-de report sheet.
$wsOverview.Range("G2").Value = "2016-09-05 11:35:28"
$wsDeDe.Range("H2").Value = "2016-09-05 11:35:28"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime for the
# first row.
$wsZhCn.Range("H2").Value = "2016-09-05 11:35:23"
$wsZhCn.Range("K2").Value = "2016-09-05 11:35:41"

# de-de: Correspond Handback DateTime for the first row.
$wsDeDe.Range("K2").Value = "2016-09-05 11:35:49"
